$wb = $excel.ActiveWorkbook

# "zh-cn" sheet: update handoff/handback datetimes for row 3
# (eddf713a-11ab-4a3b-bf75-0b113a2820e2 file) to reflect a new handback run.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-19 18:37:32"
$wsZhCn.Range("H3").Value = "2016-03-19 18:37:51"

# "de-de" sheet: update handoff/handback datetimes for row 3
# (eddf713a-11ab-4a3b-bf75-0b113a2820e2 file) to reflect a new handback run.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-19 18:37:35"
$wsDeDe.Range("H3").Value = "2016-03-19 18:37:56"
